$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8584253006459687
$ws.Range("C2").Value = 0.2883171445852273

$ws.Range("B3").Value = 0.9440780667701379
$ws.Range("C3").Value = 0.5564356437668533

$ws.Range("B4").Value = 0.2196686436106886
$ws.Range("C4").Value = 0.04393218513957309

$ws.Range("B5").Value = 0.0806535225480115
$ws.Range("C5").Value = 0.6612488798201234

$ws.Range("B6").Value = 0.903534649916815
$ws.Range("C6").Value = 0.9009844779010717
